$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 numeric updates
$ws.Range("E2").Value2 = [double]"25.00000000000047"
$ws.Range("H2").Value2 = [double]"1.128562159720617e-16"
$ws.Range("I2").Value2 = [double]"0.0308467551830367"
$ws.Range("K2").Value2 = [double]"46.3600891046285"
$ws.Range("O2").Value2 = [double]"1.540921321580579"
$ws.Range("S2").Value2 = [double]"54.07957996612102"
$ws.Range("W2").Value2 = [double]"18.86886886886922"
$ws.Range("X2").Value2 = [double]"18.51851851851887"
$ws.Range("Y2").Value2 = [double]"19.21921921921958"

# Row 2 string (CI) updates
$ws.Range("L2").Value2 = "[42.706687112084275, 50.013491097172725]"
$ws.Range("P2").Value2 = "[1.452868674633117, 1.628973968528041]"
$ws.Range("T2").Value2 = "[51.579500727986144, 56.57965920425589]"

# Row 3 numeric updates
$ws.Range("E3").Value2 = [double]"24.01000000000031"
$ws.Range("H3").Value2 = [double]"1.128562159720617e-16"
$ws.Range("K3").Value2 = [double]"48.57858145615877"
$ws.Range("O3").Value2 = [double]"0.3710790121357315"
$ws.Range("Q3").Value2 = [double]"1.798561299892754e-14"
$ws.Range("R3").Value2 = [double]"1.798561299892754e-14"
$ws.Range("S3").Value2 = [double]"49.55057264073593"
$ws.Range("W3").Value2 = [double]"22.59199199199229"
$ws.Range("X3").Value2 = [double]"22.25551551551581"
$ws.Range("Y3").Value2 = [double]"22.92846846846877"

# Row 3 string (CI) updates
$ws.Range("L3").Value2 = "[44.13285318168667, 53.02430973063086]"
$ws.Range("P3").Value2 = "[0.2830263651882694, 0.45913165908319353]"
$ws.Range("T3").Value2 = "[47.062273390142074, 52.038871891329784]"
